$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(118).Insert()

$ws.Cells.Item(118, 1).Value = 3
$ws.Cells.Item(118, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44775
$ws.Cells.Item(118, 5).Value = 5
$ws.Cells.Item(118, 6).Value = 100112039
$ws.Cells.Item(118, 7).Value = "Ciboulette"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 120
$ws.Cells.Item(118, 11).Value = 1500
$ws.Cells.Item(118, 12).Value = 1500
$ws.Cells.Item(118, 13).Value = 1500
$ws.Cells.Item(118, 14).Value = "`$/docena de atados"
$ws.Cells.Item(118, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(118, 16).Value = 500
$ws.Cells.Item(118, 17).Value = 3
$ws.Cells.Item(118, 18).Value = "Hortaliza"
